$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "20ml" and "50ml" Zithrox rows (row 5 <-> row 6):
# Row 5 becomes the 50ml item (value moves from M5 to U5)
# Row 6 becomes the 20ml item (value moves from U6 to M6)

$ws.Range("D5").Value = "Zithrox 50ml Powder for Suspension"
$ws.Range("M5").ClearContents()
$ws.Range("U5").Value = 344

$ws.Range("D6").Value = "Zithrox 20ml Powder for Suspension"
$ws.Range("U6").ClearContents()
$ws.Range("M6").Value = 3
